$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cells = @("B2", "C3", "D4", "E5", "H8", "I9", "K11", "L12", "M13", "N14", "P16", "Q17", "R18", "T20", "Y25", "Z26")

foreach ($cell in $cells) {
    $ws.Range($cell).Value = 0
}
